# "Append Data Code Drafted"
# Adds a "Season" column in front of the existing Liverpool Stats table,
# tags the existing season's row, appends a new row for the 23/24 season,
# and re-applies a column filter on the new Season column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a brand-new column before column A - this shifts the whole
#    table (and its column widths) one column to the right.
$ws.Columns("A").Insert()

# 2. Grow Table1 so it covers the new column plus a new (3rd) row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("B1:N3"))

# 3. Populate the new "Season" column.
$ws.Range("A1").Value = "Season"
$ws.Range("A2").Value = "24/25"
$ws.Range("A3").Value = "23/24"

# 4. Fill in the stats for the newly added 23/24 season row.
$ws.Range("B3").Value = "Liverpool"
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 82
$ws.Range("E3").Value = 38
$ws.Range("F3").Value = 24
$ws.Range("G3").Value = 4
$ws.Range("H3").Value = 10
$ws.Range("I3").Value = 86
$ws.Range("J3").Value = 41
$ws.Range("K3").Value = 45
$ws.Range("L3").Value = "Mo Salah"
$ws.Range("M3").Value = 5
$ws.Range("N3").Value = 65

# 5. Put a plain column filter on the new Season column (separate from the
#    Table1 autofilter), and register it as the sheet's filter database.
$ws.Range("A1:A3").AutoFilter()
$n = $ws.Names.Add("_xlnm._FilterDatabase", "='Liverpool Stats'!`$A`$1:`$A`$3")
$n.Visible = $false

# 6. Leave the cursor where the author last left it.
$ws.Range("M10").Select()
